# Applies the "healthJam forms" update to medication.xlsx:
#  - survey sheet: widen column A, change the "Current ART Regime" row to a
#    select_multiple, and add a new "Current Prophylaxis" select_multiple row
#  - choices sheet: add choice lists for the two new select_multiple questions
#  - settings sheet: picks up shared-string shifts automatically; the NOW()
#    cached value refreshes on recalculation

$wb = $excel.ActiveWorkbook

$survey  = $wb.Worksheets.Item(1)
$choices = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# survey sheet
# ---------------------------------------------------------------------------

# Widen column A (16.25 -> ~37.13 chars; Excel quantizes column widths to
# whole pixels, so 36.33 is the closest input that lands on the nearest
# reachable width)
$survey.Columns.Item(1).ColumnWidth = 36.33

# Row 50 used to be the "art_regime" text question; turn it into the new
# select_multiple question, keeping its existing formatting.
$survey.Cells.Item(50, 1).Value2 = "select_multiple tenof "
$survey.Cells.Item(50, 2).Value2 = "current_art"
$survey.Cells.Item(50, 3).Value2 = "Current ART Regime"

# Insert a brand-new row for the "Current Prophylaxis" question right after
# it (pushes the old "Notes"/"end group" rows down by one).
$survey.Rows.Item(51).Insert()

$survey.Range("A43:C43").Copy()
$survey.Range("A51:C51").PasteSpecial(-4122)
$survey.Range("D43:J43").Copy()
$survey.Range("D51:J51").PasteSpecial(-4122)
$survey.Range("K39:AB39").Copy()
$survey.Range("K51:AB51").PasteSpecial(-4122)

$survey.Cells.Item(51, 1).Value2 = "select_multiple prophylaxis or_other"
$survey.Cells.Item(51, 2).Value2 = "current_prophylaxis"
$survey.Cells.Item(51, 3).Value2 = "Current Prophylaxis"

# ---------------------------------------------------------------------------
# choices sheet
# ---------------------------------------------------------------------------

# Widen the label column (C) so the new long choice labels are readable.
$choices.Columns.Item(3).ColumnWidth = 58.1

$choicesRows = @(
    @("tenof", "tenof1", "TLD - Tenofovir/Lamivudine/Dolutegravir (1st line)"),
    @("tenof", "tenof2", "Abacvir/lamuvidine/Dolutegravir (1st line)"),
    @("tenof", "tenof3", "Zidovudine/Lamivudine +Atazanavir/Ritonavir (2nd line)"),
    @("tenof", "tenof4", "Abacavir/Lamivudine + Atazanavir/Ritonavir`n"),
    @("tenof", "tenof5", "Tenofovir/Lamivudine +Lopinavir/Ritonavir `n"),
    @("tenof", "tenof6", "Tenofovir/Lamivudine + Raltegravir (3rd line)"),
    @("tenof", "tenof7", "Tenofovir/Lamivudine/Dolutegravir+Darunavir/Ritonovir (3rd line)"),
    @("tenof", "tenof8", "Zidovudine/Lamivudine/Dolutegravir+Darunavir/Ritonavir (3rd line)"),
    @("prophylaxis", "bactrim ", "Bactrim "),
    @("prophylaxis", "azithromycin", "Azithromycin")
)

$r = 2
foreach ($row in $choicesRows) {
    $choices.Cells.Item($r, 1).Value2 = $row[0]
    $choices.Cells.Item($r, 2).Value2 = $row[1]
    $choices.Cells.Item($r, 3).Value2 = $row[2]
    $r = $r + 1
}

# Match the formatting used for the rest of the data rows in this workbook.
$survey.Range("B43").Copy()
$choices.Range("A2:C11").PasteSpecial(-4122)

$survey.Range("D25").Copy()
$choices.Range("A3").PasteSpecial(-4122)
$choices.Range("A5").PasteSpecial(-4122)
$choices.Range("A6").PasteSpecial(-4122)

# The 4th choice row (tenof4) label wraps onto a second line, so Excel
# auto-grew the row height; tenof5 also has a trailing line break in its
# label but keeps the default row height in the authored file.
$choices.Rows.Item(5).RowHeight = 30.75
$choices.Rows.Item(6).AutoFit()
